$wb = $excel.ActiveWorkbook

# --- Update the "Sampling Method Context" instruction text on the Instructions tab ---
$ws = $wb.Worksheets.Item("Instructions")

$newText = 'Enter the Context for the Sampling Method IDs that are used for sampling this parameter.  Not applicable for field measurements/observations.  If you are using the standard methods defined by MassWateR, enter the context "MassWateR".'

$ws.Range("B7").Value = $newText

# The longer instruction text now needs a taller row to show fully wrapped.
$ws.Rows.Item(7).RowHeight = 45

# --- Update the view/selection state on the Instructions tab, then restore the Meta tab as active ---
$ws.Select()
$ws.Range("B8").Select()

$wsMeta = $wb.Worksheets.Item("Meta")
$wsMeta.Select()
